$d = $word.ActiveDocument

# Pull the whole package as a WordOpenXML "mega-package" string so we can
# retarget the three <wp:docPr>/<pic:cNvPr> "name" attributes that Word's
# InlineShape object model does not expose (InlineShape has no .Name
# property over COM) while leaving everything else untouched.
$xml = $d.WordOpenXML

# -- Pearson logo in the "first page" footer (footer1.xml, docPr id="3"):
#    image2.png -> image1.png
$xml = $xml.Replace(
    '<wp:docPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="3" name="image2.png"/>',
    '<wp:docPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="3" name="image1.png"/>'
)

# -- Pearson logo in the "default" footer (footer2.xml, docPr id="2"):
#    image2.png -> image1.png
$xml = $xml.Replace(
    '<wp:docPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="2" name="image2.png"/>',
    '<wp:docPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="2" name="image1.png"/>'
)

# -- Both Pearson logo <pic:cNvPr> (id="0") nodes share identical text, so a
#    single global replace retargets both of them at once.
$xml = $xml.Replace(
    '<pic:cNvPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="0" name="image2.png"/>',
    '<pic:cNvPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="0" name="image1.png"/>'
)

# -- BTec logo in the "first page" header (header1.xml): image1.jpg -> image2.jpg
$xml = $xml.Replace(
    '<wp:docPr descr="BTec_Logo-Orange" id="1" name="image1.jpg"/>',
    '<wp:docPr descr="BTec_Logo-Orange" id="1" name="image2.jpg"/>'
)
$xml = $xml.Replace(
    '<pic:cNvPr descr="BTec_Logo-Orange" id="0" name="image1.jpg"/>',
    '<pic:cNvPr descr="BTec_Logo-Orange" id="0" name="image2.jpg"/>'
)

$d.WordOpenXML = $xml
